$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 842.25
$ws.Range("I6").Value = 123
$ws.Range("K6").Value = 369
$ws.Range("M6").Value = -257
# Row 21
$ws.Range("H21").Value = 12354
$ws.Range("I21").Value = 11976
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 11976
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = -11508
$ws.Range("N21").Value = -15936
# Row 23
$ws.Range("H23").Value = 12354
$ws.Range("I23").Value = 11976
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 11976
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = -11742
$ws.Range("N23").Value = -15468
# Row 61
$ws.Range("H61").Value = 8259.444
$ws.Range("I61").Value = 9041.875
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 27125.625
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -26953.625
$ws.Range("N61").Value = -6344
# Row 74
$ws.Range("H74").Value = 13425.363
$ws.Range("I74").Value = 8781
$ws.Range("J74").Value = 18998.6
$ws.Range("K74").Value = 8781
$ws.Range("L74").Value = 18998.6
$ws.Range("M74").Value = -7845
$ws.Range("N74").Value = -20870.6
# Row 77
$ws.Range("H77").Value = 13425.363
$ws.Range("I77").Value = 8781
$ws.Range("J77").Value = 18998.6
$ws.Range("K77").Value = 43905
$ws.Range("L77").Value = 94993
$ws.Range("M77").Value = -39225
$ws.Range("N77").Value = -104353
# Row 100
$ws.Range("H100").Value = 4678.533
$ws.Range("I100").Value = 2597
$ws.Range("J100").Value = 6499.875
$ws.Range("K100").Value = 2597
$ws.Range("L100").Value = 6499.875
$ws.Range("M100").Value = -2056
$ws.Range("N100").Value = -7581.875
# Row 115
$ws.Range("H115").Value = 1301
$ws.Range("I115").Value = 890.1111
$ws.Range("J115").Value = 4999
$ws.Range("K115").Value = 2670.3333
$ws.Range("L115").Value = 14997
$ws.Range("M115").Value = -1103.3333
$ws.Range("N115").Value = -18131

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 141.11111
$ws.Range("I4").Value = 155.625
$ws.Range("K4").Value = 155.625
$ws.Range("M4").Value = -39.625
# Row 74
$ws.Range("H74").Value = 9999
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 9999
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 32
$ws.Range("H32").Value = 950
$ws.Range("I32").Value = 950
$ws.Range("K32").Value = 950
$ws.Range("M32").Value = -566
# Row 62
$ws.Range("H62").Value = 80000
$ws.Range("J62").Value = 80000
$ws.Range("L62").Value = 80000
$ws.Range("N62").Value = -81372
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 65
$ws.Range("H65").Value = 80000
$ws.Range("J65").Value = 80000
$ws.Range("L65").Value = 240000
$ws.Range("N65").Value = -246864
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 130
$ws.Range("H130").Value = 90000
$ws.Range("J130").Value = 90000
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040

$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Range("H18").Value = 49999
$ws.Range("J18").Value = 49999
$ws.Range("L18").Value = 49999
$ws.Range("N18").Value = -50459
# Row 58
$ws.Range("H58").Value = 5638.636
$ws.Range("I58").Value = 2987.4666
$ws.Range("J58").Value = 11319.714
$ws.Range("K58").Value = 2987.4666
$ws.Range("L58").Value = 11319.714
$ws.Range("M58").Value = -2784.4666
$ws.Range("N58").Value = -11725.714
# Row 93
$ws.Range("H93").Value = 17341.4
$ws.Range("I93").Value = 8603.666999999999
$ws.Range("K93").Value = 8603.666999999999
$ws.Range("M93").Value = -6731.666999999999
# Row 132
$ws.Range("H132").Value = 3689.2856
$ws.Range("I132").Value = 3368.6365
$ws.Range("K132").Value = 10105.9095
$ws.Range("M132").Value = -7575.9095
# Row 134
$ws.Range("H134").Value = 2444.2083
$ws.Range("I134").Value = 2212.348
$ws.Range("K134").Value = 6637.044
$ws.Range("M134").Value = -4102.044
# Row 136
$ws.Range("H136").Value = 5638.636
$ws.Range("I136").Value = 2987.4666
$ws.Range("J136").Value = 11319.714
$ws.Range("K136").Value = 8962.399800000001
$ws.Range("L136").Value = 33959.142
$ws.Range("M136").Value = -6412.399800000001
$ws.Range("N136").Value = -39059.142

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7555.3335
$ws.Range("I3").Value = 7555.3335
$ws.Range("K3").Value = 22666.0005
$ws.Range("M3").Value = -22554.0005
# Row 60
$ws.Range("H60").Value = 794.2069
$ws.Range("I60").Value = 295
$ws.Range("K60").Value = 885
$ws.Range("M60").Value = -634
# Row 122
$ws.Range("H122").Value = 1132.6666
$ws.Range("J122").Value = 1199.5
$ws.Range("L122").Value = 10795.5
$ws.Range("N122").Value = -15695.5

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 19971
$ws.Range("J20").Value = 22165
$ws.Range("L20").Value = 22165
$ws.Range("N20").Value = -22655
# Row 122
$ws.Range("H122").Value = 389161.38
$ws.Range("I122").Value = 458760.1
$ws.Range("K122").Value = 1376280.3
$ws.Range("M122").Value = -1373830.3

$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 400
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 47
$ws.Range("H47").Value = 49500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 49500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 49500
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -50480
# Row 52
$ws.Range("H52").Value = 49500
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 49500
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 49500
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -49966
# Row 61
$ws.Range("H61").Value = 3059.6667
$ws.Range("I61").Value = 2200.4783
$ws.Range("K61").Value = 2200.4783
$ws.Range("M61").Value = -1998.4783
# Row 113
$ws.Range("H113").Value = 3059.6667
$ws.Range("I113").Value = 2200.4783
$ws.Range("K113").Value = 2200.4783
$ws.Range("M113").Value = -30.47830000000022

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
# Row 67
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
# Row 126
$ws.Range("H126").Value = 4393.6665
$ws.Range("I126").Value = 1989.5555
$ws.Range("K126").Value = 5968.666499999999
$ws.Range("M126").Value = -3498.666499999999
# Row 136
$ws.Range("H136").Value = 3737.4119
$ws.Range("I136").Value = 2972.7693
$ws.Range("J136").Value = 6222.5
$ws.Range("K136").Value = 8918.3079
$ws.Range("L136").Value = 18667.5
$ws.Range("M136").Value = -6368.3079
$ws.Range("N136").Value = -23767.5
